# Auto-generated edit script: updates market-price derived cells per the
# scheduled-runner diff (currentAveragePrice* / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 991
$ws.Range("I11").Value = 991
$ws.Range("K11").Value = 991
$ws.Range("M11").Value = -851
$ws.Range("H15").Value = 1550.5857
$ws.Range("I15").Value = 1550.5857
$ws.Range("K15").Value = 4651.757100000001
$ws.Range("M15").Value = -4482.757100000001
$ws.Range("H31").Value = 922145.0600000001
$ws.Range("I31").Value = 922145.0600000001
$ws.Range("K31").Value = 2766435.18
$ws.Range("M31").Value = -2766205.18
$ws.Range("H40").Value = 2114.3572
$ws.Range("I40").Value = 2163.7273
$ws.Range("J40").Value = 1933.3334
$ws.Range("K40").Value = 2163.7273
$ws.Range("L40").Value = 1933.3334
$ws.Range("M40").Value = -1988.7273
$ws.Range("N40").Value = -2283.3334
$ws.Range("H74").Value = 4391.9165
$ws.Range("I74").Value = 3957.5715
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 3957.5715
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3021.5715
$ws.Range("N74").Value = -6872
$ws.Range("H76").Value = 4466.9165
$ws.Range("I76").Value = 4700.3335
$ws.Range("J76").Value = 3766.6667
$ws.Range("K76").Value = 4700.3335
$ws.Range("L76").Value = 3766.6667
$ws.Range("M76").Value = -4385.3335
$ws.Range("N76").Value = -4396.6667
$ws.Range("H77").Value = 4391.9165
$ws.Range("I77").Value = 3957.5715
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 19787.8575
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -15107.8575
$ws.Range("N77").Value = -34360
$ws.Range("H79").Value = 4466.9165
$ws.Range("I79").Value = 4700.3335
$ws.Range("J79").Value = 3766.6667
$ws.Range("K79").Value = 4700.3335
$ws.Range("L79").Value = 3766.6667
$ws.Range("M79").Value = -3608.3335
$ws.Range("N79").Value = -5950.6667
$ws.Range("H129").Value = 929.6731
$ws.Range("I129").Value = 597.625
$ws.Range("J129").Value = 990.0454999999999
$ws.Range("K129").Value = 1792.875
$ws.Range("L129").Value = 2970.1365
$ws.Range("M129").Value = 3207.125
$ws.Range("N129").Value = -12970.1365
$ws.Range("H132").Value = 1644.2703
$ws.Range("I132").Value = 1097.6981
$ws.Range("J132").Value = 3023.7144
$ws.Range("K132").Value = 3293.0943
$ws.Range("L132").Value = 9071.143199999999
$ws.Range("M132").Value = -763.0943000000002
$ws.Range("N132").Value = -14131.1432
$ws.Range("H138").Value = 3424.0378
$ws.Range("I138").Value = 2373.182
$ws.Range("J138").Value = 4169.8066
$ws.Range("K138").Value = 7119.545999999999
$ws.Range("L138").Value = 12509.4198
$ws.Range("M138").Value = -1979.545999999999
$ws.Range("N138").Value = -22789.4198
$ws.Range("H141").Value = 6448.724
$ws.Range("I141").Value = 3007.087
$ws.Range("J141").Value = 19641.666
$ws.Range("K141").Value = 9021.261
$ws.Range("L141").Value = 58924.99800000001
$ws.Range("M141").Value = -3841.261
$ws.Range("N141").Value = -69284.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 39950
$ws.Range("J42").Value = 29900
$ws.Range("L42").Value = 29900
$ws.Range("N42").Value = -30872
$ws.Range("H61").Value = 2606.2903
$ws.Range("I61").Value = 1369.8
$ws.Range("J61").Value = 7758.3335
$ws.Range("K61").Value = 1369.8
$ws.Range("L61").Value = 7758.3335
$ws.Range("M61").Value = -1157.8
$ws.Range("N61").Value = -8182.3335
$ws.Range("H74").Value = 1402.5834
$ws.Range("I74").Value = 1829
$ws.Range("J74").Value = 805.6
$ws.Range("K74").Value = 1829
$ws.Range("L74").Value = 805.6
$ws.Range("M74").Value = -955
$ws.Range("N74").Value = -2553.6
$ws.Range("H77").Value = 1402.5834
$ws.Range("I77").Value = 1829
$ws.Range("J77").Value = 805.6
$ws.Range("K77").Value = 9145
$ws.Range("L77").Value = 4028
$ws.Range("M77").Value = -4777
$ws.Range("N77").Value = -12764
$ws.Range("H102").Value = 202024
$ws.Range("I102").Value = 2373.3333
$ws.Range("J102").Value = 501500
$ws.Range("K102").Value = 2373.3333
$ws.Range("L102").Value = 501500
$ws.Range("M102").Value = -751.3332999999998
$ws.Range("N102").Value = -504744
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 500370
$ws.Range("J105").Value = 500370
$ws.Range("L105").Value = 500370
$ws.Range("N105").Value = -507358
$ws.Range("H136").Value = 2606.2903
$ws.Range("I136").Value = 1369.8
$ws.Range("J136").Value = 7758.3335
$ws.Range("K136").Value = 4109.4
$ws.Range("L136").Value = 23275.0005
$ws.Range("M136").Value = -1559.4
$ws.Range("N136").Value = -28375.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1436.5857
$ws.Range("I31").Value = 1113.8334
$ws.Range("J31").Value = 3373.1
$ws.Range("K31").Value = 1113.8334
$ws.Range("L31").Value = 3373.1
$ws.Range("M31").Value = -818.8334
$ws.Range("N31").Value = -3963.1
$ws.Range("H34").Value = 1436.5857
$ws.Range("I34").Value = 1113.8334
$ws.Range("J34").Value = 3373.1
$ws.Range("K34").Value = 1113.8334
$ws.Range("L34").Value = 3373.1
$ws.Range("M34").Value = -911.8334
$ws.Range("N34").Value = -3777.1
$ws.Range("H35").Value = 1867.2222
$ws.Range("I35").Value = 1829.2858
$ws.Range("J35").Value = 2000
$ws.Range("K35").Value = 1829.2858
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = -1535.2858
$ws.Range("N35").Value = -2588
$ws.Range("H52").Value = 29278.166
$ws.Range("J52").Value = 29278.166
$ws.Range("L52").Value = 29278.166
$ws.Range("N52").Value = -29866.166
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H138").Value = 90190
$ws.Range("J138").Value = 90190
$ws.Range("L138").Value = 90190
$ws.Range("N138").Value = -100470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 492.82144
$ws.Range("J34").Value = 700
$ws.Range("L34").Value = 2100
$ws.Range("N34").Value = -2268
$ws.Range("H36").Value = 3373.75
$ws.Range("I36").Value = 496.66666
$ws.Range("J36").Value = 5100
$ws.Range("K36").Value = 1489.99998
$ws.Range("L36").Value = 15300
$ws.Range("M36").Value = -1320.99998
$ws.Range("N36").Value = -15638
$ws.Range("H39").Value = 2620
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2620
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7860
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -8448
$ws.Range("H68").Value = 1542.5555
$ws.Range("I68").Value = 1160
$ws.Range("J68").Value = 1733.8334
$ws.Range("K68").Value = 3480
$ws.Range("L68").Value = 5201.5002
$ws.Range("M68").Value = -2669
$ws.Range("N68").Value = -6823.5002
$ws.Range("H71").Value = 1542.5555
$ws.Range("I71").Value = 1160
$ws.Range("J71").Value = 1733.8334
$ws.Range("K71").Value = 10440
$ws.Range("L71").Value = 15604.5006
$ws.Range("M71").Value = -6384
$ws.Range("N71").Value = -23716.5006
$ws.Range("H131").Value = 992.03
$ws.Range("J131").Value = 1015.34375
$ws.Range("L131").Value = 3046.03125
$ws.Range("N131").Value = -13126.03125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2900.2727
$ws.Range("I68").Value = 1680
$ws.Range("J68").Value = 3917.1667
$ws.Range("K68").Value = 1680
$ws.Range("L68").Value = 3917.1667
$ws.Range("M68").Value = -931
$ws.Range("N68").Value = -5415.1667
$ws.Range("H71").Value = 2900.2727
$ws.Range("I71").Value = 1680
$ws.Range("J71").Value = 3917.1667
$ws.Range("K71").Value = 8400
$ws.Range("L71").Value = 19585.8335
$ws.Range("M71").Value = -4656
$ws.Range("N71").Value = -27073.8335
$ws.Range("H111").Value = 48786.8
$ws.Range("J111").Value = 48786.8
$ws.Range("L111").Value = 48786.8
$ws.Range("N111").Value = -56966.8
$ws.Range("H132").Value = 2883.6365
$ws.Range("I132").Value = 2619.04
$ws.Range("J132").Value = 3710.5
$ws.Range("K132").Value = 7857.12
$ws.Range("L132").Value = 11131.5
$ws.Range("M132").Value = -5327.12
$ws.Range("N132").Value = -16191.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63856
$ws.Range("J46").Value = 63856
$ws.Range("L46").Value = 63856
$ws.Range("N46").Value = -64318
$ws.Range("H62").Value = 4200.316
$ws.Range("I62").Value = 3980
$ws.Range("J62").Value = 4445.1113
$ws.Range("K62").Value = 3980
$ws.Range("L62").Value = 4445.1113
$ws.Range("M62").Value = -3356
$ws.Range("N62").Value = -5693.1113
$ws.Range("H65").Value = 4200.316
$ws.Range("I65").Value = 3980
$ws.Range("J65").Value = 4445.1113
$ws.Range("K65").Value = 19900
$ws.Range("L65").Value = 22225.5565
$ws.Range("M65").Value = -16780
$ws.Range("N65").Value = -28465.5565
$ws.Range("H134").Value = 63856
$ws.Range("J134").Value = 63856
$ws.Range("L134").Value = 191568
$ws.Range("N134").Value = -196638
